# Atualização automática de preços de eletricidade
# Refresh row 2 (the single data row) with the latest day's spot-price data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 46031
$ws.Range("B2").Value = 61.66
$ws.Range("C2").Value = 44.37
$ws.Range("D2").Value = 38.25
$ws.Range("E2").Value = 25.04
$ws.Range("F2").Value = 20
$ws.Range("G2").Value = 35.75
$ws.Range("H2").Value = 49.5
$ws.Range("I2").Value = 70.25
$ws.Range("J2").Value = 78.61
$ws.Range("K2").Value = 68.16
$ws.Range("L2").Value = 56.6
$ws.Range("M2").Value = 42.57
$ws.Range("N2").Value = 26.09
$ws.Range("O2").Value = 16.26
$ws.Range("P2").Value = 6.5
$ws.Range("Q2").Value = 14.99
$ws.Range("R2").Value = 48.08
$ws.Range("S2").Value = 71.28
$ws.Range("T2").Value = 85.04000000000001
$ws.Range("U2").Value = 91.25
$ws.Range("V2").Value = 87.59999999999999
$ws.Range("W2").Value = 84.26000000000001
$ws.Range("X2").Value = 83.05
$ws.Range("Y2").Value = 77.22
$ws.Range("Z2").Value = 53.43
$ws.Range("AA2").Value = "20h-24h"
$ws.Range("AB2").Value = 83.03
$ws.Range("AC2").Value = "18h-20h"
$ws.Range("AD2").Value = 88.15000000000001
$ws.Range("AE2").Value = "20h-22h"
$ws.Range("AF2").Value = 85.93000000000001
$ws.Range("AG2").Value = "1h-16h"
